# [update] => master data and use effect menu to context main page
#
# Inserts a new "seq" column before column A on the active sheet, shifting
# the existing columns (id, name, otherPurchase, price, materialPrice,
# vatPrice, gpPrice, wagePrice, platformName) one column to the right, and
# fills the new column with a 1-based sequence number for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Shift the whole sheet one column to the right to make room for "seq".
$ws.Columns.Item(1).Insert(-4161)

# Header for the new column, using the same look as the other headers.
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)
$ws.Cells.Item(1, 1).Value = "seq"

# Fill column A, rows 2..lastRow, with the row's sequence number, copying
# the data-row style from its neighboring (shifted) cell.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $r - 1
}

$excel.CutCopyMode = 0
